$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "66.338.11"
Set-TextValue "E2" "  +0.06%  "
Set-TextValue "D3" "3.563.50"
Set-TextValue "E3" "  -0.07%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "605.29"
Set-TextValue "E5" "  -0.21%  "
Set-TextValue "D6" "147.19"
Set-TextValue "E6" "  +1.76%  "
Set-TextValue "D7" "3.563.25"
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "E8" "  -0.19%  "
Set-TextValue "D9" "0.488"
Set-TextValue "E9" "  -0.02%  "
Set-TextValue "E10" "  -1.53%  "
Set-TextValue "D11" "7.86"
Set-TextValue "E11" "  +0.76%  "
Set-TextValue "D12" "0.410"
Set-TextValue "E12" "  -0.59%  "
Set-TextValue "D13" "4.167.80"
Set-TextValue "E13" "  -0.08%  "
Set-TextValue "E14" "  -1.93%  "
Set-TextValue "D15" "29.32"
Set-TextValue "E15" "  -3.34%  "
Set-TextValue "D16" "3.553.69"
Set-TextValue "E16" "  -0.21%  "
Set-TextValue "E17" "  +1.86%  "
Set-TextValue "D18" "66.310.71"
Set-TextValue "E18" "  -0.03%  "
Set-TextValue "D19" "11.08"
Set-TextValue "E19" "  -3.47%  "
Set-TextValue "D20" "6.27"
Set-TextValue "E20" "  +0.85%  "
Set-TextValue "D21" "14.81"
Set-TextValue "E21" "  -0.30%  "
Set-TextValue "D22" "422.23"
Set-TextValue "E22" "  -2.16%  "
Set-TextValue "D23" "0.607"
Set-TextValue "E23" "  -0.78%  "
Set-TextValue "D24" "77.95"
Set-TextValue "E24" "  -2.07%  "
Set-TextValue "D25" "3.703.25"
Set-TextValue "E25" "  -0.13%  "
Set-TextValue "E26" "  +0.03%  "
Set-TextValue "E27" "  -0.70%  "
Set-TextValue "D28" "9.26"
Set-TextValue "E28" "  +1.04%  "
Set-TextValue "D29" "7.95"
Set-TextValue "E29" "  -0.17%  "
Set-TextValue "D30" "2.49"
Set-TextValue "E30" "  -0.77%  "
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  -0.08%  "
Set-TextValue "D32" "3.560.21"
Set-TextValue "E32" "  +0.00%  "
Set-TextValue "E33" "  +3.07%  "
Set-TextValue "D34" "24.77"
Set-TextValue "E34" "  -2.78%  "
Set-TextValue "E35" "  +0.01%  "
Set-TextValue "D36" "1.34"
Set-TextValue "E36" "  -8.17%  "
Set-TextValue "D37" "7.67"
Set-TextValue "E37" "  -2.07%  "
Set-TextValue "E38" "  -4.13%  "
Set-TextValue "E39" "  -6.39%  "
Set-TextValue "E40" "  -0.51%  "
Set-TextValue "D41" "0.0832"
Set-TextValue "E41" "  -2.10%  "
Set-TextValue "D42" "5.14"
Set-TextValue "E42" "  -1.10%  "
Set-TextValue "D43" "0.867"
Set-TextValue "E43" "  -2.54%  "
Set-TextValue "D44" "45.79"
Set-TextValue "E44" "  -0.48%  "
Set-TextValue "D45" "1.84"
Set-TextValue "E45" "  -4.86%  "
Set-TextValue "D46" "0.999"
Set-TextValue "E46" "  -0.04%  "
Set-TextValue "D47" "2.46"
Set-TextValue "E47" "  -2.15%  "
Set-TextValue "B48" "Cosmos"
Set-TextValue "C48" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D48" "7.11"
Set-TextValue "E48" "  -0.56%  "
Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "23.16"
Set-TextValue "E49" "  -1.29%  "
Set-TextValue "E50" "  -6.21%  "
Set-TextValue "D51" "23.68"
Set-TextValue "E51" "  -5.71%  "
